$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row: year 2010 in A2 (the sheet's dimension grows to A1:G2)
$ws.Range("A2").Value = 2010

# Move the active selection to A3, mirroring where the cursor lands after
# entering the value above
$ws.Range("A3").Select()
